$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row (row 11): Right marks per question and Wrong (negative) marks per question
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Update "Total" row (row 12): total score and the "score / max" summary text
$ws.Range("B12").Value = 88
$ws.Range("E12").Value = "88 / 112"
